$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (13:27)
$ws.Range("A13:E27").EntireRow.Delete()

# Add the new "Running" activity row
$ws.Range("A12").Value = "Running"
$ws.Range("B12").Value = 22

# Re-apply the shared formula across C3:C12 so it stays a single shared group
$ws.Range("C3:C12").Formula = "=200/1/109*B3"

# Widen column A to fit the longer activity labels
$ws.Columns("A").ColumnWidth = 43 - 5/6

# Move the active selection to D1
$ws.Range("D1").Select()
